# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets to the latest scraped totals.

$wb = $excel.ActiveWorkbook

# 展览 (sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 615
$ws1.Range("F3").Value  = 10620
$ws1.Range("F8").Value  = 12654
$ws1.Range("F9").Value  = 13068
$ws1.Range("F10").Value = 1322
$ws1.Range("F12").Value = 5525
$ws1.Range("F33").Value = 170
$ws1.Range("F37").Value = 4458
$ws1.Range("F42").Value = 2272
$ws1.Range("F45").Value = 301
$ws1.Range("F46").Value = 45
$ws1.Range("F48").Value = 4308

# 演出 (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 84

# 本地生活 (sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 90

# 全部类型 (sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 615
$ws4.Range("F3").Value  = 10620
$ws4.Range("F7").Value  = 90
$ws4.Range("F8").Value  = 12654
$ws4.Range("F9").Value  = 13068
$ws4.Range("F11").Value = 1322
$ws4.Range("F13").Value = 5525
$ws4.Range("F33").Value = 170
$ws4.Range("F38").Value = 4458
$ws4.Range("F43").Value = 2272
$ws4.Range("F45").Value = 301
$ws4.Range("F46").Value = 45
$ws4.Range("F48").Value = 4308
